$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows right after the existing row 1092. This shifts every
# existing row from 1093..1156 down to 1097..1160 (dimension grows from
# A1:T1156 to A1:T1160), matching the target diff.
$ws.Rows("1093:1096").Insert()

# Fill in the 4 newly inserted rows with their data. Columns A, B, C, E, F,
# G, H, I, J, Q, R, T are constant for every "Platano" row in this sheet, so
# reuse those constants; only D (fecha), K (variedad), L (calidad), M
# (volumen), N (precio minimo), O (precio maximo), P (precio promedio
# ponderado) and S vary per row.

$newRows = @(
    @{ D = 44706; K = "Sin especificar"; L = "Maduro";          M = 820;  N = 7000; O = 7500;  P = 7232; S = 362 },
    @{ D = 44706; K = "Sin especificar"; L = "Pintón";          M = 1020; N = 7000; O = 8000;  P = 7529; S = 376 },
    @{ D = 44706; K = "Sin especificar"; L = "Primera Maduro";  M = 1000; N = 9000; O = 9500;  P = 9260; S = 463 },
    @{ D = 44706; K = "Sin especificar"; L = "Primera Pintón";  M = 2420; N = 9000; O = 10000; P = 9595; S = 480 }
)

$r = 1093
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108006
    $ws.Cells.Item($r, 10).Value = "Plátano"
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = "`$/caja 20 kilos"
    $ws.Cells.Item($r, 18).Value = "Ecuador"
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = 20
    $r = $r + 1
}
